# Fruta / hortaliza, semanal
# Reassigns the per-row Fecha (D), Volumen (M), Precio mínimo (N),
# Precio máximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# values across rows 2-13 (row 11 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value set for each row, taken (permuted) from the original data.
$newData = @{
    2  = @(44175, 25, 20000, 20000, 20000, 4000)
    3  = @(44188, 30, 15000, 15000, 15000, 3000)
    4  = @(44186, 40, 15000, 15000, 15000, 3000)
    5  = @(44931, 50, 18000, 18000, 18000, 3600)
    6  = @(44914, 56, 23000, 23000, 23000, 4600)
    7  = @(44189, 40, 15000, 15000, 15000, 3000)
    8  = @(44179, 45, 20000, 20000, 20000, 4000)
    9  = @(44907, 45, 25000, 25000, 25000, 5000)
    10 = @(44902, 35, 12000, 12000, 12000, 2400)
    12 = @(44196, 56, 15000, 15000, 15000, 3000)
    13 = @(44181, 30, 20000, 20000, 20000, 4000)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]

    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals[5]   # S - Precio $/Kg
}
